$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.007.51"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.530.02"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.24%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  +5.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "4.140.09"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.40%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "67.984.15"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "3.524.38"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "399.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.81%  "
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "24.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.879"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.42%  "
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "2.888.74"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "353.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "
